# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" text block with new rates ---
$hoja1 = $wb.Worksheets.Item("Hoja1")
$text = $hoja1.Range("A1").Value()
$text = $text -replace "7725.63", "7708.69"
$text = $text -replace "938.04", "949.09"
$hoja1.Range("A1").Value = $text

# --- tasas: update the N10/O10/N12/O12 rate figures ---
$tasas = $wb.Worksheets.Item("tasas")
$tasas.Range("N10").Value = 472
$tasas.Range("O10").Value = 3638.5
$tasas.Range("N12").Value = 3655
$tasas.Range("O12").Value = 450.002
